# Update build timestamp from "February 03 2026 17.29.55 EST" to
# "February 03 2026 18.05.36 EST" across the workbook.

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Black Eagle Coal Mine, United States, M3402, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
# Column S holds the "build_version" text in rows 2 through 10.
$newVersionText = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
for ($row = 2; $row -le 10; $row++) {
    $wsBoundaries.Cells.Item($row, 19).Value = $newVersionText  # column S = 19
}
